# Update report template: refresh the year/date stamps, left-align the
# tools table, and re-tune a few syntax-highlighting colors (ConstantTok,
# SpecialCharTok, FunctionTok, AttributeTok).

$d = $word.ActiveDocument

# --- 1. Subtitle "... - 2022" -> "... - 2024" -------------------------
$subtitle = $d.Paragraphs(2).Range
$subtitle.Find.Execute("2022", $true, $false, $false, $false, $false, $true, 1, $false, "2024", 2)

# --- 2. Date "2022-02-11" -> "2024-03-01" ------------------------------
$dateRange = $d.Paragraphs(4).Range
$dateRange.Find.Execute("2022-02-11", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-01", 2)

# --- 3. Left-align the (only) resources table --------------------------
$table = $d.Tables(1)
$table.Alignment = 0   # wdAlignRowLeft -> <w:jc w:val="left"/>

# --- 4. Syntax-highlighting style colors --------------------------------
$constantTok = $d.Styles("ConstantTok")
$constantTok.Font.Color = 153999         # 0x8f5902 (RRGGBB -> BGR)

$specialCharTok = $d.Styles("SpecialCharTok")
$specialCharTok.Font.Color = 23758       # 0xce5c00
$specialCharTok.Font.Bold = $true

$functionTok = $d.Styles("FunctionTok")
$functionTok.Font.Color = 8866336        # 0x204a87
$functionTok.Font.Bold = $true

$attributeTok = $d.Styles("AttributeTok")
$attributeTok.Font.Color = 8866336       # 0x204a87
